# "Generate Report for Archive"
#
# 1. Status text: every cell that held the shared string "Ready for handoff"
#    (Overview!E2, Overview!F2, zh-cn!C2, de-de!C2) moves to "In Translation".
# 2. The "Status" column is narrowed: Overview columns E & F (zh-cn / de-de
#    status columns) and column C ("Status") on the zh-cn and de-de sheets
#    shrink from ~17.22 chars to ~13.41 chars. Excel's ColumnWidth setter
#    snaps to a pixel grid, so we pick the input that lands on the closest
#    reachable grid point to the target width.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "In Translation" ---
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the Status columns ---
$wsOverview.Columns.Item(5).ColumnWidth = 12.43
$wsOverview.Columns.Item(6).ColumnWidth = 12.43
$wsZhCn.Columns.Item(3).ColumnWidth = 12.43
$wsDeDe.Columns.Item(3).ColumnWidth = 12.43
